# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1134059
$ws.Range("C4").Value = 3029
$ws.Range("D4").Value = 161782
$ws.Range("E4").Value = 906391
$ws.Range("G4").Value = 133
$ws.Range("H4").Value = 65886

# Row 5 - España
$ws.Range("C5").Value = 2588
$ws.Range("F5").Value = 2386

# Row 9 - Alemania
$ws.Range("B9").Value = 164197
$ws.Range("C9").Value = 120
$ws.Range("E9").Value = 28461

# Row 22 - Arabia Saudita
$ws.Range("F22").Value = 139

# Row 23 - Portugal
$ws.Range("B23").Value = 25190
$ws.Range("C23").Value = 203
$ws.Range("E23").Value = 22496

# Row 59 - Moldavia
$ws.Range("B59").Value = 4052
$ws.Range("C59").Value = 72
$ws.Range("D59").Value = 1334
$ws.Range("E59").Value = 2594
$ws.Range("G59").Value = 2
$ws.Range("H59").Value = 124

# Rows 113/114 - Maldivas moves above Mali (alphabetical reorder) with
# Maldivas getting updated figures and Mali's own figures unchanged.
$ws.Range("A113").Value = "Maldivas"
$ws.Range("B113").Value = 514
$ws.Range("C113").Value = 23
$ws.Range("D113").Value = 17
$ws.Range("E113").Value = 496
$ws.Range("F113").Value = 2
$ws.Range("H113").Value = 1

$ws.Range("A114").Value = "Mali"
$ws.Range("B114").Value = 508
$ws.Range("D114").Value = 196
$ws.Range("E114").Value = 286
$ws.Range("F114").Value = 0
$ws.Range("H114").Value = 26

# Row 117 - Jordania
$ws.Range("E117").Value = 86
$ws.Range("G117").Value = 1
$ws.Range("H117").Value = 9

# Row 152 - Suazilandia
$ws.Range("B152").Value = 108
$ws.Range("C152").Value = 2
$ws.Range("E152").Value = 95
